$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D header: "WA" -> "RD" ---
$ws.Range("D3").Value = "RD"

# --- Column widths for new G/H columns ---
$ws.Columns("G").ColumnWidth = 36.166666666666664
$ws.Columns("H").ColumnWidth = 23.166666666666664

# --- B7:B22 opcode bin values: numbers -> 4-digit binary text strings ---
$ws.Range("B7:B22").NumberFormat = "@"
$bin = @("0000","0001","0010","0011","0100","0101","0110","0111","1000","1001","1010","1011","1100","1101","1110","1111")
for ($i = 0; $i -lt 16; $i++) {
    $row = 7 + $i
    $ws.Cells.Item($row, 2).Value = $bin[$i]
}

# --- New "Examples" section in columns G/H ---
$ws.Range("G1").Value = "Examples"

$ws.Range("G2").Value = "1010_0000_0000_0001_00000000"
$ws.Range("H2").Value = "load 1 = x(0+0)"

$ws.Range("G3").Value = "1010_0000_0000_0010_00000001"
$ws.Range("H3").Value = "load 2 = x(0+1)"

$ws.Range("G4").Value = "0000_0001_0010_0011_00000000"
$ws.Range("H4").Value = "ADD 3 = R1+R2"

$ws.Range("G5").Value = "1011_0000_0011_0000_00000011"
$ws.Range("H5").Value = "STORE x(0+3) = R3"

$ws.Range("G7").Value = "1010_0000_0000_0001_00000100"
$ws.Range("H7").Value = "LOAD 1 = x(0+4)"

$ws.Range("G8").Value = "1000_0001_0000_0010_00000101"
$ws.Range("H8").Value = "ADDI R2 = R1+5"

$ws.Range("G9").Value = "1011_0000_0010_0000_00001000"
$ws.Range("H9").Value = "STORE x(0+8) = R2"

$ws.Range("G10").Value = "1100_0001_0010_0000_00000010"
$ws.Range("H10").Value = "BEQ R1==R2 +2"

$ws.Range("G11").Value = "0001_0001_0010_0011_00000000"
$ws.Range("H11").Value = "R3=R1-R2"

$ws.Range("G12").Value = "1110_0000_0000_0000_00000111"
$ws.Range("H12").Value = "JMP 7"

$ws.Range("G13").Value = "1111_0000_0000_0000_00000000"
$ws.Range("H13").Value = "HALT"

# --- Update selection to match the author's final cursor position ---
$ws.Range("F18").Select()
